$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 6).Value = 'StatQuest'
$ws.Cells.Item(3, 6).Value = 'ML Fundamentals: Cross Validation'
$ws.Cells.Item(2, 6).Value = 'A Gentle Introduction to ML'
$ws.Cells.Item(4, 6).Value = 'ML Fundamentals: Confusion Matrix'
$ws.Cells.Item(5, 6).Value = 'ML Fundamentals: Sensitivity and Specificity'
$ws.Cells.Item(6, 6).Value = 'ML Fundamentals: Bias and Variance'
$ws.Cells.Item(7, 6).Value = 'ROC and AUC, Clearly Explained!'
$ws.Cells.Item(8, 6).Value = 'ROC and AUC in R'
$ws.Cells.Item(10, 6).Value = 'SQ: Fitting a line to data, aka least squares, aka linear regression'
$ws.Cells.Item(11, 6).Value = 'SQ: Odds and Log(Odds) Clearly Explained!!!'
$ws.Cells.Item(12, 6).Value = 'SQ: Odds Ratios and Log(Odds Ratios), Clearly Explained!!!'
$ws.Cells.Item(14, 6).Value = 'SQ: Logistic Regression'
$ws.Cells.Item(15, 6).Value = 'Logistic Regression Details Pt1: Coefficients'
$ws.Cells.Item(16, 6).Value = 'Logistic Regression Details Pt2: Maximum Likelihood'
$ws.Cells.Item(17, 6).Value = 'Logistic Regression Details Pt3: R-squared and p-value'
$ws.Cells.Item(18, 6).Value = 'Satirated Models and Deviance'
$ws.Cells.Item(19, 6).Value = 'Logistic Regression in R, Clearly Explained!!!'
$ws.Cells.Item(20, 6).Value = 'Deviance Residuals'
$ws.Cells.Item(22, 6).Value = 'Regularization Part 1: Ridge Regression'
$ws.Cells.Item(23, 6).Value = 'Regularization Part 2: Lasso Regression'
$ws.Cells.Item(24, 6).Value = 'Regularization Part 3: Elastic Net Regression'
$ws.Cells.Item(25, 6).Value = 'Ridge, Lsso and Elastic-Net Regression in R'
$ws.Cells.Item(27, 6).Value = 'SQ: Principal Component Analysis (PCA), Step-by-Step'
$ws.Cells.Item(28, 6).Value = 'SQ: PCA main ideas in only 5 minutes!!!'
$ws.Cells.Item(30, 6).Value = 'SQ: PCA in R'
$ws.Cells.Item(31, 6).Value = 'SQ: PCA in Python'
$ws.Cells.Item(33, 6).Value = 'SQ: Linear Discriminant Analysis (LDA) clearly explained'
$ws.Cells.Item(34, 6).Value = 'SQ: MDS and PCoA'
$ws.Cells.Item(35, 6).Value = 'SQ: MDS and PCoA in R'
$ws.Cells.Item(37, 6).Value = 'SQ: t-SNE, Clearly Explained'
$ws.Cells.Item(38, 6).Value = 'SQ: Hierarchical Clustering'
$ws.Cells.Item(39, 6).Value = 'SQ: K-means Clustering'
$ws.Cells.Item(40, 6).Value = 'SQ: K-nearest neighbors, Clearly Explained'
$ws.Cells.Item(42, 6).Value = 'SQ: Decision Trees'
$ws.Cells.Item(43, 6).Value = 'SQ: Decision Trees, Part 2 - Feature Selection and Missing Data'
$ws.Cells.Item(45, 6).Value = 'SQ: Random Forests Part 1 - Building, Using and Evaluating'
$ws.Cells.Item(46, 6).Value = 'SQ: Random Forests Part 2 - Missing data and clustering'
$ws.Cells.Item(47, 6).Value = 'SQ: Random Forests in R'
$ws.Cells.Item(29, 6).Value = 'SQ: PCA - Practical Tips'
$ws.Cells.Item(49, 6).Value = 'Gradient Descent, Step-by-Step'
$ws.Cells.Item(50, 6).Value = 'Stochastic Gradient Descent, Clearly Explained!!!'
$ws.Cells.Item(51, 6).Value = 'AdaBoost, Clearly Explained'
$ws.Cells.Item(52, 6).Value = 'Gradient Boost Part 1: Regression Main Ideas'
$ws.Cells.Item(53, 6).Value = 'Gradient Boost Part 2: Regression Details'
$ws.Cells.Item(54, 6).Value = 'Gradient Boost Part 3: Classification'
$ws.Cells.Item(55, 6).Value = 'Gradient Boost Part 4: Classification Details'
$ws.Cells.Item(57, 6).Value = 'SQ: Fitting a curve to data, aka lowess, aka loess'
$ws.Cells.Item(58, 6).Value = 'Statistics Fundamentals: Population Parameters'
$ws.Cells.Item(59, 6).Value = 'Pricipal Component Analysis (PCA) clearly explained (2015)'
